{"js": "// \"Change tube code square\" \u2014 resize/reposition the floating ${tube_code}\n// square (a wordprocessingShape textbox anchored in the document body).\n// The shape keeps its right edge fixed while growing leftwards: the left\n// offset decreases and the width increases by the same EMU amount.\n\nconst EMU_PER_PT = 12700;\n\n// Target absolute position/size, expressed in EMU (taken from the anchor's\n// <wp:posOffset>/<wp:extent> after the edit), then converted to points \u2014\n// the unit Word's object model (Shape.Left / Shape.Width) uses.\nconst NEW_LEFT_EMU = 3633175;\nconst NEW_WIDTH_EMU = 2846999;\n\nconst shapes = context.document.body.shapes;\nshapes.load(\"items\");\nawait context.sync();\n\nif (shapes.items.length === 0) {\n  throw new Error(\"No floating shape found in the document body.\");\n}\n\n// There is a single floating shape in this layout \u2014 the tube-code square.\nconst shape = shapes.items[0];\nshape.load(\"left,top,width,height,name\");\nawait context.sync();\n\n// Only the horizontal position and width change; top/height stay as-is.\nshape.left = NEW_LEFT_EMU / EMU_PER_PT;\nshape.width = NEW_WIDTH_EMU / EMU_PER_PT;\n\nawait context.sync();\n", "ps1": "# \"Change tube code square\" \u2014 resize/reposition the floating ${tube_code}\n# square (a wordprocessingShape textbox anchored in the document body).\n# The shape keeps its right edge fixed while growing leftwards: the left\n# offset decreases and the width increases by the same EMU amount.\n\n$EMU_PER_PT = 12700\n\n# Target absolute position/size, expressed in EMU (taken from the anchor's\n# <wp:posOffset>/<wp:extent> after the edit), then converted to points \u2014\n# the unit Word's object model (Shape.Left / Shape.Width) uses.\n$NEW_LEFT_EMU = 3633175\n$NEW_WIDTH_EMU = 2846999\n\n$d = $word.ActiveDocument\n\nif ($d.Shapes.Count -eq 0) {\n    throw \"No floating shape found in the document body.\"\n}\n\n# There is a single floating shape in this layout \u2014 the tube-code square.\n$shp = $d.Shapes.Item(1)\n\n# Only the horizontal position and width change; top/height stay as-is.\n$shp.Left = $NEW_LEFT_EMU / $EMU_PER_PT\n$shp.Width = $NEW_WIDTH_EMU / $EMU_PER_PT\n"}
